# Fix units on household emissions.
# Update the header labels on the "CO2 emissions" sheet so that the
# "[tonnes]" unit suffix becomes "[million tonnes]". Dependent sheets
# that pull these headers via formulas will pick up the change
# automatically on recalculation.

$wb = $excel.ActiveWorkbook

$coSheet = $wb.Worksheets.Item("CO2 emissions")
$coSheet.Range("B1").Value = "CO2 emiss [million tonnes]"
$coSheet.Range("C1").Value = "CO2 emiss elect alloc[million tonnes]"
$coSheet.Range("D1").Value = "CO2 emissions allocated [million tonnes]"

# Autofit column B on the "data to plot" sheet so the longer label fits.
$plotSheet = $wb.Worksheets.Item("data to plot")
$plotSheet.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

$excel.Calculate()

# Leave the selections where the editing session ended up.
$coSheet.Range("D1").Select() | Out-Null

$elecSheet = $wb.Worksheets.Item("US elect cons 2018")
$elecSheet.Range("B1").Select() | Out-Null

$plotSheet.Range("C12").Select() | Out-Null

